$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# --- Row 2: Lakkalikoori / 1 / Success ---
$ws.Rows.Item(2).RowHeight = 29
$ws.Range("A2").Value = "Lakkalikoori"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Success"
$ws.Range("D2").Clear()

# --- Row 3: Chang / 1 / Success ---
$ws.Rows.Item(3).RowHeight = 43.5
$ws.Range("A3").Value = "Chang"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Success"

# --- Row 4: Ipoh Coffee / 1 / Failed / Unable to find beveraged named 'Ipoh Coffee' ---
$ws.Rows.Item(4).RowHeight = 29
$ws.Range("A4").Value = "Ipoh Coffee"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Failed"
$ws.Range("D4").Value = "Unable to find beveraged named 'Ipoh Coffee'"
$ws.Range("D4").WrapText = $true

# --- Row 5: Laughing Lumberjack Lager / 1 / Success ---
$ws.Rows.Item(5).RowHeight = 43.5
$ws.Range("A5").Value = "Laughing Lumberjack Lager"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "Success"

# --- Rows 6, 7, 9: fully cleared (content + formatting), so the <row> disappears ---
$ws.Range("A6:D6").Clear()
$ws.Range("A7:D7").Clear()
$ws.Range("A9:D9").Clear()

# --- Rows 8 and 10: keep row height but clear all cell content/format ---
$ws.Rows.Item(8).RowHeight = 29
$ws.Rows.Item(10).RowHeight = 43.5
$ws.Range("A8:D8").Clear()
$ws.Range("A10:D10").Clear()

# --- Column widths: B and C become a fixed width, matching column A's style family ---
$ws.Range("B1:C1").ColumnWidth = 8.3
